$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.285.62"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").Value = "1.873.38"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5058"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3932"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.488"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.32%  "

$ws.Range("D14").Value = "1.878.02"
$ws.Range("E14").Value = "  +4.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.429"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.63%  "

$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06604"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.10%  "

$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.36%  "

$ws.Range("D23").Value = "28.340.56"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("E25").Value = "  +2.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.565"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.25%  "

$ws.Range("D27").Value = "2.094.91"
$ws.Range("E27").Value = "  +4.17%  "

$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1066"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.066"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06741"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.500"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02400"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.16%  "

$ws.Range("E38").Value = "  +2.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("E40").Value = "  +2.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.183"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.61%  "

$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.17%  "

$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.271"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("E50").Value = "  +1.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06854"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.73%  "
